$wb = $excel.ActiveWorkbook

# Add the new "correlation" worksheet after Sheet1 and make it the active sheet.
$sheet1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "correlation"

# Fill in the correlation matrix data.
$ws.Range("B1").Value = "all"
$ws.Range("C1").Value = "replace"
$ws.Range("D1").Value = "delete_insert"
$ws.Range("E1").Value = "pos<50%"

$ws.Range("A2").Value = "replace"
$ws.Range("B2").Value = 0.97

$ws.Range("A3").Value = "delete_insert"
$ws.Range("B3").Value = 0.59
$ws.Range("C3").Value = 0.37

$ws.Range("A4").Value = "pos<50%"
$ws.Range("B4").Value = 0.94
$ws.Range("C4").Value = 0.92
$ws.Range("D4").Value = 0.51

$ws.Range("A5").Value = "pos>=50%"
$ws.Range("B5").Value = 0.78
$ws.Range("C5").Value = 0.74
$ws.Range("D5").Value = 0.52
$ws.Range("E5").Value = 0.52

# Match the recorded selection on the new sheet.
$result = $ws.Range("G12").Select()
